# Actualización automática del tracker
# Rellena las celdas de resultado (G) y profit (H) para las filas 14 y 16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G14").Value = "Acierto"
$ws.Range("H14").Value = 1.63

$ws.Range("G16").Value = "Acierto"
$ws.Range("H16").Value = 0.83
